$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update Version value (row 3, column B): 0.1.1 -> 0.2.0
$ws.Range("B3").Value = "0.2.0"

# Update Date value (row 8, column B)
$ws.Range("B8").Value = "2023-10-19T17:05:12+00:00"

# Insert a new row after the "Contact" row (row 10) for the new
# "Jurisdiction" / "iso:code:3166:FR" property, pushing Description,
# Purpose, Copyright and Immutable down by one row.
$ws.Rows.Item(11).Insert()

# The freshly inserted row doesn't carry the surrounding data-row
# formatting, so copy it over from the row above ("Contact") before
# filling in the new values.
$ws.Range("A10:B10").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0

$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = "iso:code:3166:FR"
